$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.410.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.076.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.373"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.604.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.481.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.066.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.500"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0908"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0674"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.110.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.669"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.286.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0254"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.936"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0877"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "250.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.59%  "
